$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates
$ws.Range("B1").Value = "Modèle"
$ws.Range("C1").Value = "Score"

# Row 2 updates
$ws.Range("B2").Value = "Factuality and Readability"
$ws.Range("C2").Value = "0.8414592146873474 / 2"
$ws.Range("D2").Value = 0.8958436548709869
$ws.Range("E2").Value = 2

# Row 3 updates
$ws.Range("B3").Value = "Factuality and Readability"
$ws.Range("C3").Value = "0.9502280950546265 / 2"
